$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 26 ("Scyliorhinus canicula") - this shifts rows 27:60 up to 26:59
$ws.Rows.Item(26).Delete()

# After the shift, update H (Numb) cells that are 0 to -1, for the species/discard rows
# that originally had Numb = 0 (now located at these rows after the deletion).
$zeroRows = @(36, 37, 38, 42, 52, 54, 56, 57)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 8).Value = -1
}

# Update the I (RF) column value from 61.35244444444444 to 61.40244444444446
# for rows 30 through 59 (rows that carry the running RF total after the deleted row).
for ($r = 30; $r -le 59; $r++) {
    $ws.Cells.Item($r, 9).Value = 61.40244444444446
}

